$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.050673333333334
$ws.Range("H2").Value = 9.15202
$ws.Range("I2").Value = 0.03589373762206555
$ws.Range("J2").Value = 0.03589373762206555
$ws.Range("M2").Value = 8.432170666666666
$ws.Range("N2").Value = 25.296512
$ws.Range("O2").Value = 0.0153412147997323
$ws.Range("P2").Value = 0.01534121479973231
$ws.Range("Q2").Value = 25.72379819491556
$ws.Range("R2").Value = 231.51418375424
$ws.Range("S2").Value = 0.0005506535388253402
$ws.Range("T2").Value = 0.0005506535388253403
# Row 3
$ws.Range("G3").Value = 3.050673333333334
$ws.Range("H3").Value = 9.15202
$ws.Range("I3").Value = 0.03589373762206555
$ws.Range("J3").Value = 0.03589373762206555
$ws.Range("M3").Value = 211.5004576666667
$ws.Range("N3").Value = 634.5013730000001
$ws.Range("O3").Value = 0.3847969970689267
$ws.Range("P3").Value = 0.3847969970689267
$ws.Range("Q3").Value = 645.2188061914957
$ws.Range("R3").Value = 5806.96925572346
$ws.Range("S3").Value = 0.01381180245055078
$ws.Range("T3").Value = 0.01381180245055078
# Row 4
$ws.Range("G4").Value = 3.050673333333334
$ws.Range("H4").Value = 9.15202
$ws.Range("I4").Value = 0.03589373762206555
$ws.Range("J4").Value = 0.03589373762206555
$ws.Range("M4").Value = 149.6042426666667
$ws.Range("N4").Value = 448.812728
$ws.Range("O4").Value = 0.2721850532240109
$ws.Range("P4").Value = 0.2721850532240109
$ws.Range("Q4").Value = 456.3936736567289
$ws.Range("R4").Value = 4107.54306291056
$ws.Range("S4").Value = 0.009769738885070593
$ws.Range("T4").Value = 0.009769738885070595
# Row 5
$ws.Range("G5").Value = 3.050673333333334
$ws.Range("H5").Value = 9.15202
$ws.Range("I5").Value = 0.03589373762206555
$ws.Range("J5").Value = 0.03589373762206555
$ws.Range("M5").Value = 180.1047823333333
$ws.Range("N5").Value = 540.314347
$ws.Range("O5").Value = 0.3276767349073302
$ws.Range("P5").Value = 0.3276767349073302
$ws.Range("Q5").Value = 549.4408566701045
$ws.Range("R5").Value = 4944.96771003094
$ws.Range("S5").Value = 0.01176154274761884
$ws.Range("T5").Value = 0.01176154274761884
# Row 6
$ws.Range("I6").Value = 0.7605628985450701
$ws.Range("J6").Value = 0.7605628985450702
$ws.Range("M6").Value = 8.432170666666666
$ws.Range("N6").Value = 25.296512
$ws.Range("O6").Value = 0.0153412147997323
$ws.Range("P6").Value = 0.01534121479973231
$ws.Range("Q6").Value = 545.069079255936
$ws.Range("R6").Value = 4905.621713303424
$ws.Range("S6").Value = 0.01166795879528693
$ws.Range("T6").Value = 0.01166795879528693
# Row 7
$ws.Range("I7").Value = 0.7605628985450701
$ws.Range("J7").Value = 0.7605628985450702
$ws.Range("M7").Value = 211.5004576666667
$ws.Range("N7").Value = 634.5013730000001
$ws.Range("O7").Value = 0.3847969970689267
$ws.Range("P7").Value = 0.3847969970689267
$ws.Range("S7").Value = 0.2926623194421817
$ws.Range("T7").Value = 0.2926623194421817
# Row 8
$ws.Range("I8").Value = 0.7605628985450701
$ws.Range("J8").Value = 0.7605628985450702
$ws.Range("M8").Value = 149.6042426666667
$ws.Range("N8").Value = 448.812728
$ws.Range("O8").Value = 0.2721850532240109
$ws.Range("P8").Value = 0.2721850532240109
$ws.Range("Q8").Value = 9670.658959199784
$ws.Range("R8").Value = 87035.93063279804
$ws.Range("S8").Value = 0.2070138530206979
$ws.Range("T8").Value = 0.2070138530206979
# Row 9
$ws.Range("I9").Value = 0.7605628985450701
$ws.Range("J9").Value = 0.7605628985450702
$ws.Range("M9").Value = 180.1047823333333
$ws.Range("N9").Value = 540.314347
$ws.Range("O9").Value = 0.3276767349073302
$ws.Range("P9").Value = 0.3276767349073302
$ws.Range("Q9").Value = 11642.26291862144
$ws.Range("R9").Value = 104780.366267593
$ws.Range("S9").Value = 0.2492187672869036
$ws.Range("T9").Value = 0.2492187672869036
# Row 10
$ws.Range("G10").Value = 16.398149
$ws.Range("H10").Value = 49.194447
$ws.Range("I10").Value = 0.1929380151136699
$ws.Range("J10").Value = 0.19293801511367
$ws.Range("M10").Value = 8.432170666666666
$ws.Range("N10").Value = 25.296512
$ws.Range("O10").Value = 0.0153412147997323
$ws.Range("P10").Value = 0.01534121479973231
$ws.Range("Q10").Value = 138.2719909854293
$ws.Range("R10").Value = 1244.447918868864
$ws.Range("S10").Value = 0.002959903532892808
$ws.Range("T10").Value = 0.002959903532892809
# Row 11
$ws.Range("G11").Value = 16.398149
$ws.Range("H11").Value = 49.194447
$ws.Range("I11").Value = 0.1929380151136699
$ws.Range("J11").Value = 0.19293801511367
$ws.Range("M11").Value = 211.5004576666667
$ws.Range("N11").Value = 634.5013730000001
$ws.Range("O11").Value = 0.3847969970689267
$ws.Range("P11").Value = 0.3847969970689267
$ws.Range("Q11").Value = 3468.216018386192
$ws.Range("R11").Value = 31213.94416547573
$ws.Range("S11").Value = 0.07424196883617938
$ws.Range("T11").Value = 0.07424196883617938
# Row 12
$ws.Range("G12").Value = 16.398149
$ws.Range("H12").Value = 49.194447
$ws.Range("I12").Value = 0.1929380151136699
$ws.Range("J12").Value = 0.19293801511367
$ws.Range("M12").Value = 149.6042426666667
$ws.Range("N12").Value = 448.812728
$ws.Range("O12").Value = 0.2721850532240109
$ws.Range("P12").Value = 0.2721850532240109
$ws.Range("Q12").Value = 2453.232662280157
$ws.Range("R12").Value = 22079.09396052142
$ws.Range("S12").Value = 0.05251484391264926
$ws.Range("T12").Value = 0.05251484391264927
# Row 13
$ws.Range("G13").Value = 16.398149
$ws.Range("H13").Value = 49.194447
$ws.Range("I13").Value = 0.1929380151136699
$ws.Range("J13").Value = 0.19293801511367
$ws.Range("M13").Value = 180.1047823333333
$ws.Range("N13").Value = 540.314347
$ws.Range("O13").Value = 0.3276767349073302
$ws.Range("P13").Value = 0.3276767349073302
$ws.Range("Q13").Value = 2953.385056314567
$ws.Range("R13").Value = 26580.46550683111
$ws.Range("S13").Value = 0.06322129883194848
$ws.Range("T13").Value = 0.06322129883194849
# Row 14
$ws.Range("G14").Value = 0.9013676666666667
$ws.Range("H14").Value = 2.704103
$ws.Range("I14").Value = 0.01060534871919427
$ws.Range("J14").Value = 0.01060534871919427
$ws.Range("M14").Value = 8.432170666666666
$ws.Range("N14").Value = 25.296512
$ws.Range("O14").Value = 0.0153412147997323
$ws.Range("P14").Value = 0.01534121479973231
$ws.Range("Q14").Value = 7.600485998748444
$ws.Range("R14").Value = 68.404373988736
$ws.Range("S14").Value = 0.0001626989327272251
$ws.Range("T14").Value = 0.0001626989327272252
# Row 15
$ws.Range("G15").Value = 0.9013676666666667
$ws.Range("H15").Value = 2.704103
$ws.Range("I15").Value = 0.01060534871919427
$ws.Range("J15").Value = 0.01060534871919427
$ws.Range("M15").Value = 211.5004576666667
$ws.Range("N15").Value = 634.5013730000001
$ws.Range("O15").Value = 0.3847969970689267
$ws.Range("P15").Value = 0.3847969970689267
$ws.Range("Q15").Value = 190.6396740259354
$ws.Range("R15").Value = 1715.757066233419
$ws.Range("S15").Value = 0.004080906340014741
$ws.Range("T15").Value = 0.004080906340014743
# Row 16
$ws.Range("G16").Value = 0.9013676666666667
$ws.Range("H16").Value = 2.704103
$ws.Range("I16").Value = 0.01060534871919427
$ws.Range("J16").Value = 0.01060534871919427
$ws.Range("M16").Value = 149.6042426666667
$ws.Range("N16").Value = 448.812728
$ws.Range("O16").Value = 0.2721850532240109
$ws.Range("P16").Value = 0.2721850532240109
$ws.Range("Q16").Value = 134.8484271358871
$ws.Range("R16").Value = 1213.635844222984
$ws.Range("S16").Value = 0.002886617405593087
$ws.Range("T16").Value = 0.002886617405593088
# Row 17
$ws.Range("G17").Value = 0.9013676666666667
$ws.Range("H17").Value = 2.704103
$ws.Range("I17").Value = 0.01060534871919427
$ws.Range("J17").Value = 0.01060534871919427
$ws.Range("M17").Value = 180.1047823333333
$ws.Range("N17").Value = 540.314347
$ws.Range("O17").Value = 0.3276767349073302
$ws.Range("P17").Value = 0.3276767349073302
$ws.Range("Q17").Value = 162.3406274073046
$ws.Range("R17").Value = 1461.065646665741
$ws.Range("S17").Value = 0.003475126040859213
$ws.Range("T17").Value = 0.003475126040859214

Write-Output "Updated 182 cells with new TPM values"